$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "Staff_1"
$ws.Range("B2").Value = "M1"
$ws.Range("C2").Value = "M1"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = "M1"
$ws.Range("F2").Value = "M1"
$ws.Range("G2").Value = "A1"
$ws.Range("H2").Value = "M3"
$ws.Range("I2").Value = "M1"
$ws.Range("J2").Value = "M3"
$ws.Range("K2").Value = "M1"
$ws.Range("L2").Value = "M1"
$ws.Range("M2").Value = "M1"
$ws.Range("N2").Value = "A1"
$ws.Range("O2").Value = "M1"
$ws.Range("P2").Value = "M1"
$ws.Range("Q2").Value = "M1"
$ws.Range("R2").Value = "M1"
$ws.Range("S2").Value = "M1"
$ws.Range("T2").Value = "M3"
$ws.Range("U2").Value = "A1"
$ws.Range("V2").Value = "M1"
$ws.Range("W2").Value = "M1"
$ws.Range("X2").Value = "M1"
$ws.Range("Y2").Value = "M3"
$ws.Range("Z2").Value = "M1"
$ws.Range("AA2").Value = "M1"
$ws.Range("AB2").Value = "A1"
$ws.Range("AC2").Value = "M1"
$ws.Range("A3").Value = "Staff_2"
$ws.Range("B3").Value = "M2"
$ws.Range("C3").Value = "M3"
$ws.Range("D3").Value = "M1"
$ws.Range("E3").Value = "M1"
$ws.Range("F3").Value = "M2"
$ws.Range("G3").Value = "A1"
$ws.Range("H3").Value = "M1"
$ws.Range("I3").Value = "M1"
$ws.Range("J3").Value = "M2"
$ws.Range("K3").Value = "M2"
$ws.Range("L3").Value = "M2"
$ws.Range("M3").Value = "M1"
$ws.Range("N3").Value = "A2"
$ws.Range("O3").Value = "M3"
$ws.Range("P3").Value = "M3"
$ws.Range("Q3").Value = "M2"
$ws.Range("R3").Value = "M1"
$ws.Range("S3").Value = "M2"
$ws.Range("T3").Value = "M1"
$ws.Range("U3").Value = "A2"
$ws.Range("V3").Value = "M1"
$ws.Range("W3").Value = "M1"
$ws.Range("X3").Value = "M1"
$ws.Range("Y3").Value = "M2"
$ws.Range("Z3").Value = "M3"
$ws.Range("AA3").Value = "M1"
$ws.Range("AB3").Value = "A1"
$ws.Range("AC3").Value = "M1"
$ws.Range("A4").Value = "Staff_3"
$ws.Range("B4").Value = "M1"
$ws.Range("C4").Value = "M1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = "M1"
$ws.Range("F4").Value = "M1"
$ws.Range("G4").Value = "A1"
$ws.Range("H4").Value = "M3"
$ws.Range("I4").Value = "M1"
$ws.Range("J4").Value = "M1"
$ws.Range("K4").Value = "M1"
$ws.Range("L4").Value = "M1"
$ws.Range("M4").Value = "M3"
$ws.Range("N4").Value = "A1"
$ws.Range("O4").Value = "M1"
$ws.Range("P4").Value = "M1"
$ws.Range("Q4").Value = "M1"
$ws.Range("R4").Value = "M1"
$ws.Range("S4").Value = "M3"
$ws.Range("T4").Value = "M1"
$ws.Range("U4").Value = "A1"
$ws.Range("V4").Value = "M1"
$ws.Range("W4").Value = "M1"
$ws.Range("X4").Value = "M1"
$ws.Range("Y4").Value = "M1"
$ws.Range("Z4").Value = "M1"
$ws.Range("AA4").Value = "M3"
$ws.Range("AB4").Value = "A1"
$ws.Range("AC4").Value = "M1"
$ws.Range("A5").Value = "Staff_4"
$ws.Range("B5").Value = "M2"
$ws.Range("C5").Value = "M3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = "M2"
$ws.Range("F5").Value = "M1"
$ws.Range("G5").Value = "A1"
$ws.Range("H5").Value = "M1"
$ws.Range("I5").Value = "M1"
$ws.Range("J5").Value = "M2"
$ws.Range("K5").Value = "M2"
$ws.Range("L5").Value = "M3"
$ws.Range("M5").Value = "M2"
$ws.Range("N5").Value = "A2"
$ws.Range("O5").Value = "M2"
$ws.Range("P5").Value = "M1"
$ws.Range("Q5").Value = "M1"
$ws.Range("R5").Value = "M1"
$ws.Range("S5").Value = "M1"
$ws.Range("T5").Value = "M3"
$ws.Range("U5").Value = "A1"
$ws.Range("V5").Value = "M1"
$ws.Range("W5").Value = "M1"
$ws.Range("X5").Value = "M1"
$ws.Range("Y5").Value = "M2"
$ws.Range("Z5").Value = "M3"
$ws.Range("AA5").Value = "M1"
$ws.Range("AB5").Value = "A1"
$ws.Range("AC5").Value = "M2"
$ws.Range("A6").Value = "Staff_5"
$ws.Range("B6").Value = "M1"
$ws.Range("C6").Value = "M2"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = "M2"
$ws.Range("F6").Value = "M1"
$ws.Range("G6").Value = "A1"
$ws.Range("H6").Value = "M3"
$ws.Range("I6").Value = "M3"
$ws.Range("J6").Value = "M2"
$ws.Range("K6").Value = "M2"
$ws.Range("L6").Value = "M2"
$ws.Range("M6").Value = "M2"
$ws.Range("N6").Value = "A2"
$ws.Range("O6").Value = "M1"
$ws.Range("P6").Value = "M1"
$ws.Range("Q6").Value = "M2"
$ws.Range("R6").Value = "M2"
$ws.Range("S6").Value = "M1"
$ws.Range("T6").Value = "M3"
$ws.Range("U6").Value = "A1"
$ws.Range("V6").Value = "M3"
$ws.Range("W6").Value = "M1"
$ws.Range("X6").Value = "M2"
$ws.Range("Y6").Value = "M1"
$ws.Range("Z6").Value = "M3"
$ws.Range("AA6").Value = "M2"
$ws.Range("AB6").Value = "A2"
$ws.Range("AC6").Value = "M2"
$ws.Range("A7").Value = "Staff_6"
$ws.Range("B7").Value = "M1"
$ws.Range("C7").Value = "M1"
$ws.Range("D7").Value = "M1"
$ws.Range("E7").Value = "M1"
$ws.Range("F7").Value = "M1"
$ws.Range("G7").Value = "A1"
$ws.Range("H7").Value = "M3"
$ws.Range("I7").Value = "M3"
$ws.Range("J7").Value = "M1"
$ws.Range("K7").Value = "M1"
$ws.Range("L7").Value = "M1"
$ws.Range("M7").Value = "M1"
$ws.Range("N7").Value = "A1"
$ws.Range("O7").Value = "M1"
$ws.Range("P7").Value = "M1"
$ws.Range("Q7").Value = "M1"
$ws.Range("R7").Value = "M1"
$ws.Range("S7").Value = "M3"
$ws.Range("T7").Value = "M1"
$ws.Range("U7").Value = "A1"
$ws.Range("V7").Value = "M1"
$ws.Range("W7").Value = "M1"
$ws.Range("X7").Value = "M1"
$ws.Range("Y7").Value = "M1"
$ws.Range("Z7").Value = "M1"
$ws.Range("AA7").Value = "M3"
$ws.Range("AB7").Value = "A1"
$ws.Range("AC7").Value = "M1"
$ws.Range("A8").Value = "Staff_7"
$ws.Range("B8").Value = "M1"
$ws.Range("C8").Value = "M1"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = "M1"
$ws.Range("F8").Value = "M3"
$ws.Range("G8").Value = "A2"
$ws.Range("H8").Value = "M2"
$ws.Range("I8").Value = "M2"
$ws.Range("J8").Value = "M2"
$ws.Range("K8").Value = "M1"
$ws.Range("L8").Value = "M1"
$ws.Range("M8").Value = "M2"
$ws.Range("N8").Value = "A2"
$ws.Range("O8").Value = "M3"
$ws.Range("P8").Value = "M1"
$ws.Range("Q8").Value = "M2"
$ws.Range("R8").Value = "M1"
$ws.Range("S8").Value = "M1"
$ws.Range("T8").Value = "M3"
$ws.Range("U8").Value = "A1"
$ws.Range("V8").Value = "M3"
$ws.Range("W8").Value = "M2"
$ws.Range("X8").Value = "M1"
$ws.Range("Y8").Value = "M1"
$ws.Range("Z8").Value = "M2"
$ws.Range("AA8").Value = "M3"
$ws.Range("AB8").Value = "A1"
$ws.Range("AC8").Value = "M2"
$ws.Range("A9").Value = "Staff_8"
$ws.Range("B9").Value = "M1"
$ws.Range("C9").Value = "M2"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = "M2"
$ws.Range("F9").Value = "M2"
$ws.Range("G9").Value = "A2"
$ws.Range("H9").Value = "M3"
$ws.Range("I9").Value = "M2"
$ws.Range("J9").Value = "M2"
$ws.Range("K9").Value = "M1"
$ws.Range("L9").Value = "M3"
$ws.Range("M9").Value = "M1"
$ws.Range("N9").Value = "A2"
$ws.Range("O9").Value = "M2"
$ws.Range("P9").Value = "M1"
$ws.Range("Q9").Value = "M1"
$ws.Range("R9").Value = "M2"
$ws.Range("S9").Value = "M1"
$ws.Range("T9").Value = "M3"
$ws.Range("U9").Value = "A1"
$ws.Range("V9").Value = "M3"
$ws.Range("W9").Value = "M1"
$ws.Range("X9").Value = "M2"
$ws.Range("Y9").Value = "M1"
$ws.Range("Z9").Value = "M1"
$ws.Range("AA9").Value = "M2"
$ws.Range("AB9").Value = "A1"
$ws.Range("AC9").Value = "M3"
$ws.Range("A10").Value = "Staff_9"
$ws.Range("B10").Value = "M1"
$ws.Range("C10").Value = "M2"
$ws.Range("D10").Value = "M3"
$ws.Range("E10").Value = "M1"
$ws.Range("F10").Value = "M2"
$ws.Range("G10").Value = "A2"
$ws.Range("H10").Value = "M2"
$ws.Range("I10").Value = "M2"
$ws.Range("J10").Value = "M1"
$ws.Range("K10").Value = "M2"
$ws.Range("L10").Value = "M1"
$ws.Range("M10").Value = "M3"
$ws.Range("N10").Value = "A1"
$ws.Range("O10").Value = "M2"
$ws.Range("P10").Value = "M1"
$ws.Range("Q10").Value = "M1"
$ws.Range("R10").Value = "M3"
$ws.Range("S10").Value = "M1"
$ws.Range("T10").Value = "M1"
$ws.Range("U10").Value = "A1"
$ws.Range("V10").Value = "M2"
$ws.Range("W10").Value = "M1"
$ws.Range("X10").Value = "M2"
$ws.Range("Y10").Value = "M1"
$ws.Range("Z10").Value = "M2"
$ws.Range("AA10").Value = "M2"
$ws.Range("AB10").Value = "A1"
$ws.Range("AC10").Value = "M3"